$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stressor1")

# --- D column: convert numeric adcap_score values to text labels ---
$ws.Range("D2").Value = "none"
$ws.Range("D3").Value = "low"
$ws.Range("D4").Value = "medium"
$ws.Range("D5").Value = "high"
$ws.Range("D6").Value = "low"
$ws.Range("D7").Value = "medium"
$ws.Range("D8").Value = "high"
$ws.Range("D9").Value = "high"
$ws.Range("D10").Value = "none"
$ws.Range("D11").Value = "none"
$ws.Range("D12").Value = "low"
$ws.Range("D13").Value = "low"
$ws.Range("D14").Value = "low"
$ws.Range("D15").Value = "medium"
$ws.Range("D16").Value = "medium"
$ws.Range("D17").Value = "medium"
$ws.Range("D18").Value = "high"
$ws.Range("D19").Value = "high"
$ws.Range("D20").Value = "none"
$ws.Range("D21").Value = "low"
$ws.Range("D22").Value = "low"
$ws.Range("D23").Value = "medium"
$ws.Range("D24").Value = "medium"
$ws.Range("D25").Value = "high"
$ws.Range("D26").Value = "none"
$ws.Range("D27").Value = "low"
$ws.Range("D28").Value = "low"
$ws.Range("D29").Value = "medium"
$ws.Range("D30").Value = "high"
$ws.Range("D31").Value = "none"
$ws.Range("D32").Value = "low"
$ws.Range("D33").Value = "low"
$ws.Range("D34").Value = "medium"
$ws.Range("D35").Value = "medium"
$ws.Range("D36").Value = "high"
$ws.Range("D37").Value = "high"
$ws.Range("D38").Value = "none"
$ws.Range("D39").Value = "high"
$ws.Range("D40").Value = "low"
$ws.Range("D41").Value = "none"
$ws.Range("D42").Value = "low"
$ws.Range("D43").Value = "medium"
$ws.Range("D44").Value = "medium"
$ws.Range("D45").Value = "high"
$ws.Range("D46").Value = "none"
$ws.Range("D47").Value = "low"
$ws.Range("D48").Value = "medium"
$ws.Range("D49").Value = "high"
$ws.Range("D50").Value = "high"
$ws.Range("D51").Value = "high"
$ws.Range("D52").Value = "none"
$ws.Range("D53").Value = "none"
$ws.Range("D54").Value = "low"
$ws.Range("D55").Value = "high"
$ws.Range("D56").Value = "high"
$ws.Range("D63").Value = "none"
$ws.Range("D64").Value = "high"
$ws.Range("D65").Value = "none"
$ws.Range("D66").Value = "low"
$ws.Range("D67").Value = "high"
$ws.Range("D68").Value = "NA"
$ws.Range("D75").Value = "none"
$ws.Range("D76").Value = "high"

# --- E7: new note referencing literature review ---
$ws.Range("E7").Value = "lit review"

# --- view/window state ---
$excel.ActiveWindow.Zoom = 125
$ws.Range("D77").Select()

# best-effort: reposition the workbook window (xWindow/yWindow)
try {
    $wb.Windows.Item(1).Left = 2120
    $wb.Windows.Item(1).Top = 120
} catch {
    # not all hosts expose window positioning; ignore
}
